# Rename the embedded logo pictures (wp:docPr / pic:cNvPr "name" attribute)
# in the document's header and footers:
#   - Pearson logo (footers)   : image2.png -> image1.png
#   - BTec logo    (header)    : image1.jpg -> image2.jpg
#
# The logos live in the section's Headers/Footers, not in the body, so they
# are reached via Section.Headers / Section.Footers rather than
# ActiveDocument.InlineShapes.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # --- Footers: Pearson Edexcel logo -> rename image2.png to image1.png ---
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.Name -eq "image2.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }

    # --- Headers: BTEC logo -> rename image1.jpg to image2.jpg ---
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.Name -eq "image1.jpg") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
}
